$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Another Person" (row 4) no longer has a phone number - clear C4 (keeps its number-as-text style)
$ws.Range("C4").Value = $null

# Row 5 becomes a new contact: "Third Person"
$ws.Range("A5").Value = "Third Person"
$ws.Range("B5").Value = "email@gmail.com"
$ws.Range("C5").Value = "+91 202934880"

# Row 6 (new): "Forth Person" (re-uses the old "hisemail@gmail.com" contact email)
$ws.Range("A6").Value = "Forth Person"
$ws.Range("B6").Value = "hisemail@gmail.com"
$ws.Range("C6").Value = "+91 238907 234"

# Row 7 (new): "Fifth person"
$ws.Range("A7").Value = "Fifth person"
$ws.Range("B7").Value = "heremail@gmail.com"
$ws.Range("C7").Value = "+1 238974293"

# Rebuild every mailto hyperlink in the Email Address column (B3:B7) so the old B5 link
# ("hisemail@gmail.com") doesn't linger after that cell's content changed.
$ws.Range("B3").Hyperlinks.Delete() | Out-Null

$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:youremail@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:personemail@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:email@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:hisemail@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:heremail@gmail.com") | Out-Null

# Hyperlinks.Add applies its own formatting - reapply the workbook's "Hyperlink" cell style
# so B3:B7 all match
$ws.Range("B3:B7").Style = "Hyperlink"

# Give C5:C7 the same "phone number" text style used by C3/C4
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Match the workbook's recorded selection after the edit
$ws.Range("B8").Select() | Out-Null
